$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update existing rows 218-222 with revised figures (columns B..O only -
# P, Q, R are unchanged for these rows per the source data update).
# ---------------------------------------------------------------------------

# Row 218
$ws.Range("B218").Value = 3962
$ws.Range("C218").Value = 884
$ws.Range("D218").Value = -55
$ws.Range("E218").Value = 719
$ws.Range("F218").Value = 220
$ws.Range("G218").Value = 3688
$ws.Range("I218").Value = 3545
$ws.Range("J218").Value = 3635
$ws.Range("L218").Value = -1138
$ws.Range("M218").Value = 527
$ws.Range("N218").Value = 110
$ws.Range("O218").Value = 174

# Row 219
$ws.Range("B219").Value = 5563
$ws.Range("C219").Value = 6178
$ws.Range("D219").Value = 6026
$ws.Range("E219").Value = -106
$ws.Range("F219").Value = 258
$ws.Range("G219").Value = 629
$ws.Range("I219").Value = 660
$ws.Range("J219").Value = 26
$ws.Range("N219").Value = 306
$ws.Range("O219").Value = -321

# Row 220
$ws.Range("B220").Value = 2442
$ws.Range("C220").Value = 1769
$ws.Range("D220").Value = 264
$ws.Range("E220").Value = 1249
$ws.Range("F220").Value = 257
$ws.Range("G220").Value = -529
$ws.Range("I220").Value = -612
$ws.Range("J220").Value = -586
$ws.Range("M220").Value = 2057
$ws.Range("N220").Value = 771
$ws.Range("O220").Value = 1239

# Row 221
$ws.Range("B221").Value = 5035
$ws.Range("C221").Value = 3654
$ws.Range("E221").Value = 1559
$ws.Range("F221").Value = 786
$ws.Range("I221").Value = 2400
$ws.Range("M221").Value = -394
$ws.Range("N221").Value = -71
$ws.Range("O221").Value = -260

# Row 222
$ws.Range("B222").Value = 1754
$ws.Range("C222").Value = -116
$ws.Range("D222").Value = 410
$ws.Range("E222").Value = -305
$ws.Range("F222").Value = -222
$ws.Range("G222").Value = 2427
$ws.Range("I222").Value = 2206
$ws.Range("J222").Value = 2448
$ws.Range("M222").Value = 272
$ws.Range("N222").Value = 386
$ws.Range("O222").Value = -169

# ---------------------------------------------------------------------------
# Append the new monthly row 223 ("01-06-2021").
#
# A223 must hold the text "01-06-2021" as a shared string, not get
# auto-converted to a date serial number. Writing it as a text formula and
# then collapsing the formula to its static value via copy / paste-special
# (values) keeps the cell a plain text string cell with no extra number
# format / style being minted.
# ---------------------------------------------------------------------------
$ws.Range("A223").Formula = "=""01-06-2021"""
$ws.Range("A223").Copy() | Out-Null
$ws.Range("A223").PasteSpecial(-4163) | Out-Null

$ws.Range("B223").Value = 497
$ws.Range("C223").Value = 1346
$ws.Range("D223").Value = 1251
$ws.Range("E223").Value = 1125
$ws.Range("F223").Value = -1030
$ws.Range("G223").Value = 1128
$ws.Range("H223").Value = 187
$ws.Range("I223").Value = 941
$ws.Range("J223").Value = 1199
$ws.Range("K223").Value = -258
$ws.Range("L223").Value = -949
$ws.Range("M223").Value = -1029
$ws.Range("N223").Value = 11
$ws.Range("O223").Value = -1028
$ws.Range("P223").Value = -11
$ws.Range("Q223").Value = 0
$ws.Range("R223").Value = 0
